# Insert a new weekly price row for "Agrícola del Norte S.A. de Arica - Caigua".
# The new record is inserted right after the existing row 56, pushing every
# row below it down by one (old row 57 -> new row 58, ..., old row 89 -> new row 90).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 57; everything from the old row 57
# downward shifts down by one row (old 89 becomes new 90).
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new observation.
$ws.Cells.Item(57, 1).Value = 1
$ws.Cells.Item(57, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(57, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(57, 4).Value = 44596
$ws.Cells.Item(57, 5).Value = 15
$ws.Cells.Item(57, 6).Value = 100112036
$ws.Cells.Item(57, 7).Value = "Caigua"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 130
$ws.Cells.Item(57, 11).Value = 8000
$ws.Cells.Item(57, 12).Value = 9000
$ws.Cells.Item(57, 13).Value = 8500
$ws.Cells.Item(57, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(57, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(57, 16).Value = 425
$ws.Cells.Item(57, 17).Value = 20
$ws.Cells.Item(57, 18).Value = "Hortaliza"
